$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.713.82"
$ws.Range("E2").Value = "  -0.94%  "
$ws.Range("D3").Value = "3.787.44"
$ws.Range("E3").Value = "  +1.16%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.82"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.10"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.01%  "
$ws.Range("D7").Value = "3.786.92"
$ws.Range("E7").Value = "  +1.16%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.160"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.52%  "
$ws.Range("E11").Value = "  -1.55%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.449"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.56%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000253"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.72%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.04"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.50%  "
$ws.Range("D15").Value = "4.419.89"
$ws.Range("E15").Value = "  +1.14%  "
$ws.Range("D16").Value = "3.763.75"
$ws.Range("E16").Value = "  +0.09%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.50"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +3.43%  "
$ws.Range("D18").Value = "67.662.91"
$ws.Range("E18").Value = "  -0.91%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.04"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.32%  "
$ws.Range("E20").Value = "  -0.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.01"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -6.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "460.30"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.65%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.696"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.62%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000153"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +6.17%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.30"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.58%  "
$ws.Range("E26").Value = "  +1.74%  "
$ws.Range("E27").Value = "  -0.85%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.04"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.55%  "
$ws.Range("E29").Value = "  +0.14%  "
$ws.Range("B30").Value = "WrappedeETH"
$ws.Range("C30").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D30").Value = "3.931.40"
$ws.Range("E30").Value = "  +1.08%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.78"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.93%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.24"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +4.48%  "
$ws.Range("B33").Value = "NEARProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.22"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.62%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.65"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.33%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.11"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.12%  "
$ws.Range("B36").Value = "Binance-PegBSC-USD"
$ws.Range("C36").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.05%  "
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.100"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.34%  "
$ws.Range("B38").Value = "dogwifhat"
$ws.Range("C38").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.34"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.43%  "
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.138"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.59%  "
$ws.Range("B40").Value = "Mantle"
$ws.Range("C40").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.996"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.16%  "
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.76"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.28%  "
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.10%  "
$ws.Range("B43").Value = "Arweave"
$ws.Range("C43").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "46.17"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +7.06%  "
$ws.Range("B44").Value = "USDe"
$ws.Range("C44").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "48.18"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +3.52%  "
$ws.Range("B46").Value = "TheGraph"
$ws.Range("C46").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.299"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.08%  "
$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "149.41"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +3.27%  "
$ws.Range("B48").Value = "Cosmos"
$ws.Range("C48").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.33"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.33%  "
$ws.Range("B49").Value = "Bittensor"
$ws.Range("C49").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "394.27"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +2.07%  "
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.82"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -4.09%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "26.64"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +4.33%  "
